$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell="D2"; Value="302.97"},
    @{Cell="E2"; Value="5.64%"},
    @{Cell="D3"; Value="31.72"},
    @{Cell="E3"; Value="7.70%"},
    @{Cell="D4"; Value="5.213"},
    @{Cell="E4"; Value="2.16%"},
    @{Cell="D5"; Value="0.07269"},
    @{Cell="E5"; Value="7.86%"},
    @{Cell="D6"; Value="7.793"},
    @{Cell="E6"; Value="5.98%"},
    @{Cell="D7"; Value="3.743"},
    @{Cell="E7"; Value="8.68%"},
    @{Cell="D8"; Value="1.460"},
    @{Cell="E8"; Value="5.36%"},
    @{Cell="D9"; Value="0.9045"},
    @{Cell="E9"; Value="-1.31%"},
    @{Cell="D10"; Value="0.01655"},
    @{Cell="E10"; Value="2,460.62%"},
    @{Cell="D11"; Value="0.1671"},
    @{Cell="E11"; Value="5.13%"},
    @{Cell="D12"; Value="0.07455"},
    @{Cell="E12"; Value="8.97%"},
    @{Cell="D13"; Value="0.07939"},
    @{Cell="E13"; Value="3.51%"},
    @{Cell="D14"; Value="0.02975"},
    @{Cell="E14"; Value="2.05%"},
    @{Cell="D15"; Value="0.09942"},
    @{Cell="E15"; Value="10.67%"},
    @{Cell="D16"; Value="0.001500"},
    @{Cell="E16"; Value="-5.31%"},
    @{Cell="D17"; Value="0.04542"},
    @{Cell="E17"; Value="1.66%"},
    @{Cell="D18"; Value="0.006458"},
    @{Cell="E18"; Value="3.08%"},
    @{Cell="D19"; Value="3.470"},
    @{Cell="E19"; Value="0.50%"},
    @{Cell="D20"; Value="2.227"},
    @{Cell="E20"; Value="-0.09%"},
    @{Cell="D21"; Value="0.3333"},
    @{Cell="E21"; Value="4.22%"},
    @{Cell="D22"; Value="0.1330"},
    @{Cell="E22"; Value="1.88%"},
    @{Cell="D23"; Value="4.298"},
    @{Cell="E23"; Value="6.01%"},
    @{Cell="D24"; Value="0.1632"},
    @{Cell="E24"; Value="3.21%"},
    @{Cell="D25"; Value="0.001224"},
    @{Cell="E25"; Value="2.61%"},
    @{Cell="D26"; Value="0.004411"},
    @{Cell="E26"; Value="6.61%"},
    @{Cell="E27"; Value="8.92%"},
    @{Cell="D28"; Value="0.0001748"},
    @{Cell="E28"; Value="8.18%"},
    @{Cell="D40"; Value="0.04485"},
    @{Cell="E40"; Value="5.23%"},
    @{Cell="D41"; Value="0.007196"},
    @{Cell="E41"; Value="5.81%"},
    @{Cell="D42"; Value="0.1343"},
    @{Cell="E42"; Value="8.10%"},
    @{Cell="D43"; Value="0.002339"},
    @{Cell="E43"; Value="5.05%"},
    @{Cell="D44"; Value="0.01340"},
    @{Cell="E44"; Value="12.03%"},
    @{Cell="D45"; Value="0.00006073"},
    @{Cell="E45"; Value="6.66%"},
    @{Cell="E46"; Value="-3.59%"},
    @{Cell="D47"; Value="0.01613"},
    @{Cell="E47"; Value="7.21%"},
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
